# RAP3 workbook cleanup — preparing for install of RAP3 on Azure
#
# The sharedStrings table contained two unused PascalCase duplicates
# ("OrgAbbrName" / "OrgFullName") that were only referenced by the
# "Identity Provider data" sheet's helper/legend rows (B9/C9), plus a
# stray duplicate "Organization" label in E13. Clearing those cells
# drops the now-unused strings from the shared-string table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("E13").Value = ""

# Restore the last active selection to D21 (as last left by the author).
$ws.Range("D21").Select()
